# Re-ran program query after series finished
#
# The "Art on the Farm" event only had a "No" (not attended) tally row.
# Re-running the program query produced updated capacity/availability
# figures for that row and added the matching "Yes" (attended) tally row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Locate the existing "Art on the Farm" / "No" row.
$used = $ws.UsedRange
$rowCount = $used.Rows.Count
$targetRow = 0
for ($r = 2; $r -le $rowCount; $r++) {
    $name = $ws.Cells.Item($r, 1).Value2
    $attended = $ws.Cells.Item($r, 7).Value2
    if ($name -eq "Art on the Farm" -and $attended -eq "No") {
        $targetRow = $r
        break
    }
}

$newRow = $targetRow + 1

# Insert a new row right after it; Excel copies the formatting of the row above.
$ws.Rows.Item($newRow).Insert()

# Update the existing "No" row with the re-queried numbers.
$ws.Cells.Item($targetRow, 5).Value2 = 55
$ws.Cells.Item($targetRow, 6).Value2 = 1
$ws.Cells.Item($targetRow, 8).Value2 = 45

# Populate the newly inserted "Yes" row with the same event info.
$ws.Cells.Item($newRow, 1).Value2 = $ws.Cells.Item($targetRow, 1).Value2
$ws.Cells.Item($newRow, 2).Value2 = $ws.Cells.Item($targetRow, 2).Value2
$ws.Cells.Item($newRow, 3).Value2 = $ws.Cells.Item($targetRow, 3).Value2
$ws.Cells.Item($newRow, 4).Value2 = $ws.Cells.Item($targetRow, 4).Value2
$ws.Cells.Item($newRow, 5).Value2 = 55
$ws.Cells.Item($newRow, 6).Value2 = 1
$ws.Cells.Item($newRow, 7).Value2 = "Yes"
$ws.Cells.Item($newRow, 8).Value2 = 11
